$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first data row (2007年) is being retired from this series; drop it
# entirely so the remaining years (2010/2012/2015/2017) shift up one row,
# matching the refreshed data dump.
$ws.Rows("2:2").Delete()
